$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type (one cell at a time -- this COM layer only applies
# NumberFormat to the first area of a multi-area Range) for price cells
# whose new values would otherwise be auto-parsed as numbers, so they
# stay text, matching the original inline-string cell type.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values per the diff.
$ws.Range("D2").Value = '26.985.27'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.824.85'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").Value = '311.69'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '0.4626'
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").Value = '0.3702'
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("D9").Value = '0.07340'
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Value = '0.8748'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").Value = '0.07916'
$ws.Range("E11").Value = '  +4.08%  '
$ws.Range("D12").Value = '19.78'
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("D13").Value = '1.882.17'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '5.334'
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").Value = '6.538'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '91.25'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '0.000008866'
$ws.Range("E18").Value = '  +2.50%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '27.018.15'
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").Value = '5.102'
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").Value = '10.54'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '2.051.48'
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("D25").Value = '153.16'
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").Value = '2.039'
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").Value = '5.131'
$ws.Range("D30").Value = '115.52'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = '0.08882'
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").Value = '2.966'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = '0.7279'
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = '4.436'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").Value = '2.469'
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.071'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01946'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").Value = '0.05222'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = '2.946'
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").Value = '7.106'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").Value = '0.5155'
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").Value = '0.1622'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D44").Value = '8.173'
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").Value = '0.4832'
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '102.76'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '1.633'
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").Value = '64.83'
$ws.Range("E51").Value = '  +0.53%  '
